$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "275.64"
Set-TextValue $ws.Range("E2") "-0.92%"
Set-TextValue $ws.Range("D3") "27.37"
Set-TextValue $ws.Range("E3") "1.68%"
Set-TextValue $ws.Range("D4") "4.809"
Set-TextValue $ws.Range("E4") "-2.22%"
Set-TextValue $ws.Range("D5") "0.06350"
Set-TextValue $ws.Range("E5") "-0.73%"
Set-TextValue $ws.Range("D6") "6.954"
Set-TextValue $ws.Range("E6") "-0.58%"
Set-TextValue $ws.Range("D7") "1.339"
Set-TextValue $ws.Range("E7") "10.76%"
Set-TextValue $ws.Range("D8") "0.8792"
Set-TextValue $ws.Range("E8") "-0.93%"
Set-TextValue $ws.Range("E9") "1.94%"
Set-TextValue $ws.Range("D10") "0.05060"
Set-TextValue $ws.Range("E10") "-3.58%"
Set-TextValue $ws.Range("D11") "0.07509"
Set-TextValue $ws.Range("E11") "1.24%"
Set-TextValue $ws.Range("D12") "0.02970"
Set-TextValue $ws.Range("E12") "-4.73%"
Set-TextValue $ws.Range("D13") "0.09031"
Set-TextValue $ws.Range("E13") "-0.36%"
Set-TextValue $ws.Range("D14") "0.001568"
Set-TextValue $ws.Range("E14") "-0.03%"
Set-TextValue $ws.Range("D15") "0.0006422"
Set-TextValue $ws.Range("E15") "1.32%"
Set-TextValue $ws.Range("D16") "0.005690"
Set-TextValue $ws.Range("E16") "-5.76%"
Set-TextValue $ws.Range("D17") "3.449"
Set-TextValue $ws.Range("E17") "-1.20%"
Set-TextValue $ws.Range("D18") "3.304"
Set-TextValue $ws.Range("E18") "-1.60%"
Set-TextValue $ws.Range("E19") "0.16%"
Set-TextValue $ws.Range("D21") "0.1354"
Set-TextValue $ws.Range("E21") "1.68%"
Set-TextValue $ws.Range("D22") "3.908"
Set-TextValue $ws.Range("E22") "-0.78%"
Set-TextValue $ws.Range("D23") "0.04403"
Set-TextValue $ws.Range("E23") "1.36%"
Set-TextValue $ws.Range("E24") "-0.79%"
Set-TextValue $ws.Range("D25") "0.003863"
Set-TextValue $ws.Range("E25") "5.01%"
Set-TextValue $ws.Range("E26") "-0.22%"
Set-TextValue $ws.Range("D27") "0.0001937"
Set-TextValue $ws.Range("E27") "13.93%"
Set-TextValue $ws.Range("D40") "0.04176"
Set-TextValue $ws.Range("E40") "2.61%"
Set-TextValue $ws.Range("D41") "0.006855"
Set-TextValue $ws.Range("E41") "3.14%"
Set-TextValue $ws.Range("E42") "0.44%"
Set-TextValue $ws.Range("D43") "0.002029"
Set-TextValue $ws.Range("E43") "-14.17%"
Set-TextValue $ws.Range("D44") "0.01154"
Set-TextValue $ws.Range("E44") "-10.26%"
Set-TextValue $ws.Range("D45") "0.00005163"
Set-TextValue $ws.Range("E45") "-2.01%"
Set-TextValue $ws.Range("D46") "1.487"
Set-TextValue $ws.Range("E46") "-36.88%"
Set-TextValue $ws.Range("D47") "0.02300"
Set-TextValue $ws.Range("E47") "8.36%"
